# double weighting issue fixed
$wb = $excel.ActiveWorkbook

# --- Final Rankings sheet: insert a "Rank" column before "WEC Design" ---
$wsRankings = $wb.Worksheets.Item("Final Rankings")
$wsRankings.Columns("A").Insert()
$wsRankings.Columns("D:D").Delete()

$wsRankings.Range("A1").Value = "Rank"
$wsRankings.Range("B1").Value = "WEC Design"
$wsRankings.Range("C1").Value = "Closeness to Ideal"

$wsRankings.Range("A2").Value = 1
$wsRankings.Range("B2").Value = "Oscillating Surge Flap"
$wsRankings.Range("C2").Value = 0.6769679823150073

$wsRankings.Range("A3").Value = 2
$wsRankings.Range("B3").Value = "Oscillating Water Column"
$wsRankings.Range("C3").Value = 0.3505119622780461

$wsRankings.Range("A4").Value = 3
$wsRankings.Range("B4").Value = "Point Absorber"
$wsRankings.Range("C4").Value = 0.2044203284084603

# --- Final AHP Theme Weights sheet: corrected weighting values ---
$wsWeights = $wb.Worksheets.Item("Final AHP Theme Weights")

$wsWeights.Range("B2").Value = 0.2211016960397099
$wsWeights.Range("C2").Value = 0.2047978497170779
$wsWeights.Range("D2").Value = 0.1843855192755456
$wsWeights.Range("E2").Value = 0.1896748365578238
$wsWeights.Range("F2").Value = 0.2000400984098427

$wsWeights.Range("B3").Value = 0.2343673884268225
$wsWeights.Range("C3").Value = 0.2172606660340768
$wsWeights.Range("D3").Value = 0.1603201453220851
$wsWeights.Range("E3").Value = 0.3110541138651856
$wsWeights.Range("F3").Value = 0.07699768635183006

$wsWeights.Range("B4").Value = 0.1587790613253449
$wsWeights.Range("C4").Value = 0.1831033408056535
$wsWeights.Range("D4").Value = 0.2418516638648966
$wsWeights.Range("E4").Value = 0.1183195610630221
$wsWeights.Range("F4").Value = 0.2979463729410828
